# netCrypto.xlsx update (matches "Add files via upload" diff):
#  - SheetName1!T2 value changes from 631912 to 633500
#  - The active selection on SheetName1 moves from T3 to T2

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SheetName1")

# Make sure we're on the right sheet before touching the selection.
$ws.Activate()

# Update the USD Amount figure in T2.
$ws.Range("T2").Value = 633500

# Move the active cell / selection to T2 (was T3 before the edit).
$ws.Range("T2").Select()
